# Applies the "new version with timestamp" update:
#  - inserts 3 new low-stock product rows (ARBATEG, ATOMOXAPEX, RISPADEX)
#  - refreshes the grand-total and the generated-at timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $range as a genuine text value (shared-string),
# while preserving whatever number format / style is already on the cell
# (mirrors how the source report stores numeric-looking text like "1" or
# "35.0000" as text rather than as a number).
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("Z1")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# Helper: insert one brand-new product row at row $rowNum, cloning the
# layout (row height, per-cell styles, merged cells) from the row that is
# currently at $templateRow (BEFORE the insert pushes it down), then fill
# in the product's data.
function Insert-ProductRow($rowNum, $templateRow, $seq, $name, $balance, $limit, $price, $sellPrice, $count) {
    $ws.Rows("$rowNum`:$rowNum").Insert()

    $templateRange = $ws.Range("A$templateRow`:Q$templateRow")
    $newRange = $ws.Range("A$rowNum`:Q$rowNum")

    $templateRange.Copy()
    $newRange.PasteSpecial(-4122)  # xlPasteFormats (styles + row look)
    $ws.Rows("$rowNum`:$rowNum").RowHeight = $ws.Rows("$templateRow`:$templateRow").RowHeight

    $ws.Range("A$rowNum`:B$rowNum").Merge()
    $ws.Range("C$rowNum`:G$rowNum").Merge()
    $ws.Range("H$rowNum`:K$rowNum").Merge()
    $ws.Range("L$rowNum`:M$rowNum").Merge()
    $ws.Range("N$rowNum`:O$rowNum").Merge()

    $ws.Range("A$rowNum").Value = $seq
    Set-TextValue $ws.Range("C$rowNum") $name
    Set-TextValue $ws.Range("H$rowNum") $balance
    Set-TextValue $ws.Range("L$rowNum") $limit
    Set-TextValue $ws.Range("N$rowNum") $price
    Set-TextValue $ws.Range("P$rowNum") $sellPrice
    Set-TextValue $ws.Range("Q$rowNum") $count
}

# 1) ARBATEG 2% SUSP. 100ML -- new row 7, pushes old BETADERM (row 7) down
Insert-ProductRow 7 8 1 "ARBATEG 2% SUSP. 100ML" "0:0" "1" "35.00" "35.0000" "1:0"

# 2) ATOMOXAPEX 18MG 30 CAPS. -- new row 8, pushes BETADERM (now row 8) down
Insert-ProductRow 8 9 2 "ATOMOXAPEX 18MG 30 CAPS." "2:2" "1" "120.00" "120.0000" "1:0"

# 3) RISPADEX 1MG/ML SYRUP 100ML -- new row 22, pushes VIDROP (now row 22) down
Insert-ProductRow 22 22 16 "RISPADEX 1MG/ML SYRUP 100ML" "0:0" "1" "33.00" "33.0000" "1:0"

# Renumber the sequence column ("م") for every product row, now rows 7..25
for ($r = 7; $r -le 25; $r++) {
    $ws.Range("A$r").Value = $r - 6
}

$ws.Range("Z1").Clear()

# Grand total (P26, was P23) picks up the 3 new selling prices
$ws.Range("P26").Value = 951.5

# Refresh the "generated at" timestamp footer (A27, was A24)
Set-TextValue $ws.Range("A27") "Thursday, 14 August, 2025 10:54 AM"
$ws.Range("Z1").Clear()
